$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark "Completed?" column (G) for the rows that now have blocks assigned:
# Row 19 - C_add -> "Yes"
# Row 20 - C_rb  -> "yes"
# Row 22 - g     -> "Yes"
$ws.Range("G19").Value = "Yes"
$ws.Range("G20").Value = "yes"
$ws.Range("G22").Value = "Yes"

# Update selection to match the final active cell in the diff
$ws.Range("G22").Select()
